$d = $word.ActiveDocument

# --- 1. Title paragraph: merge "Anja & David Wedding - Ceremony Scr" + "ipt"
#        into a single run with fixed text "...Ceremony Script" ---
$p1 = $d.Paragraphs.Item(1)
$titleRange = $d.Range($p1.Range.Start, $p1.Range.End)
$titleRange.Text = "Anja & David Wedding " + [char]0x2014 + " Ceremony Script"

# --- 2. Date paragraph: change style to TextBody, right align, change text ---
$p2 = $d.Paragraphs.Item(2)
$dateRange = $d.Range($p2.Range.Start, $p2.Range.End)
$dateRange.Text = "2022-07-29 10:30am"
$p2.Style = "TextBody"
$p2.Range.ParagraphFormat.Alignment = 2

# --- 3. Insert a new, empty TextBody paragraph right after the date paragraph ---
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$p3.Style = "TextBody"

# --- 4. Fix the "elcome!" typo -> "Welcome!" ---
$d.Content.Find.Execute("elcome! Thank you so much for being here.", $false, $false, $false, $false, $false, $true, 1, $false, "Welcome! Thank you so much for being here.", 2) | Out-Null

# --- 5. Convert all italic runs to bold throughout the document ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Font.Italic = $true
$find.Replacement.ClearFormatting()
$find.Replacement.Font.Italic = $false
$find.Replacement.Font.Bold = $true
$find.Text = ""
$find.Replacement.Text = ""
$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null
